$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row text updates
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "60TB ZFS backup 2023.0502"
$ws.Range("B1").Value = "Price/ea"
$ws.Range("D1").Value = "Subtotal P*Q"

# ---------------------------------------------------------------------------
# 2. New "Alt" column (F)
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "Alt"
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Columns("F").ColumnWidth = 11.53

# ---------------------------------------------------------------------------
# 3. Widen column A
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 26.12

# ---------------------------------------------------------------------------
# 4. Insert a new row above the old "PC Power supply, 8x SATA" row (row 7),
#    which pushes that row (and everything below it) down by one.
# ---------------------------------------------------------------------------
$ws.Rows(7).Insert()

# New row 7 becomes the "PC power supply, 6x SATA" entry (the newly preferred
# part); it inherits row 7's old bold+yellow formatting from the insert.
$ws.Range("A7").Value = "PC power supply, 6x SATA"
$ws.Range("B7").Value = 40
$ws.Range("C7").Value = 1
$ws.Range("D7").Formula = "=B7*C7"
$ws.Range("E7").Value = "https://www.amazon.com/Enermax-Cyberbron-ECB500AWT-Non-Modular-Warranty/dp/B08K1ZBYPZ"
$ws.Range("E7").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("E7"), "https://www.amazon.com/Enermax-Cyberbron-ECB500AWT-Non-Modular-Warranty/dp/B08K1ZBYPZ", "", "", "https://www.amazon.com/Enermax-Cyberbron-ECB500AWT-Non-Modular-Warranty/dp/B08K1ZBYPZ")
$ws.Range("F7").Value = 40

# Row 8 is now the old "PC Power supply, 8x SATA" row -- demoted to an
# alternate: clear its price/qty (kept bold, but drop the yellow highlight)
# and record its old price in the new Alt column.
$ws.Range("B8").ClearContents()
$ws.Range("B8").Interior.ColorIndex = -4142
$ws.Range("C8").ClearContents()
$ws.Range("D8").Formula = "=B8*C8"
$ws.Range("F8").Value = 70

# ---------------------------------------------------------------------------
# 5. Grand Total row formatting (now row 14 after the insert)
# ---------------------------------------------------------------------------
$ws.Range("A14").Font.Bold = $true
$ws.Range("A14").HorizontalAlignment = -4152
$ws.Range("D14").NumberFormat = "[$$-409]#,##0.00;[RED]\-[$$-409]#,##0.00"

# ---------------------------------------------------------------------------
# 6. Reposition / resize the picture to match the shifted layout
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Top = 179.2063
$shp.Left = 337.011
$shp.Width = 632.296062992126
$shp.Height = 456.2929133858268
